# New weekly price record for "Hortaliza, Macroferia Regional de Talca - Pepino dulce".
# The new record is inserted as row 14 (pushing the former rows 14-21 down to 15-22),
# matching the data the rest of the sheet is sorted/grouped by (date ascending per quality).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14; existing rows 14-21 shift down to 15-22.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with this week's record.
$ws.Cells.Item(14, 1).Value  = 5
$ws.Cells.Item(14, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(14, 3).Value  = 'Maule'
$ws.Cells.Item(14, 4).Value  = 44596
$ws.Cells.Item(14, 5).Value  = 7
$ws.Cells.Item(14, 6).Value  = 100112043
$ws.Cells.Item(14, 7).Value  = 'Pepino dulce'
$ws.Cells.Item(14, 8).Value  = 'Cultivar IV Región'
$ws.Cells.Item(14, 9).Value  = 'Primera'
$ws.Cells.Item(14, 10).Value = 150
$ws.Cells.Item(14, 11).Value = 14000
$ws.Cells.Item(14, 12).Value = 14000
$ws.Cells.Item(14, 13).Value = 14000
$ws.Cells.Item(14, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(14, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(14, 16).Value = 778
$ws.Cells.Item(14, 17).Value = 18
$ws.Cells.Item(14, 18).Value = 'Hortaliza'

# Give the new date cell the same date/time number format used by the other "Fecha" cells.
$ws.Cells.Item(14, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
